# Added New Mac-Address and Document Types
# Append 5 new rows to the master-reg_center_device sheet:
#   regcntr_id 10002, device_id 3000176..3000180,
#   lang_code "eng", is_active TRUE, cr_by "superadmin", cr_dtimes "now()"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 157
$regcntrId = 10002
$startDeviceId = 3000176
$rowCount = 5

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $regcntrId
    $ws.Cells.Item($r, 2).Value = $startDeviceId + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Match the author's final scroll position / selection in the sheet view
$excel.ActiveWindow.ScrollRow = 152
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C158").Select()
